# Auto-generated Excel COM-interop script
# Applies numeric value updates to specific cells across multiple sheets
# as described by the upstream OOXML diff (scheduled data refresh).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1600.6316
$ws.Range("J17").Value = 1719.1875
$ws.Range("L17").Value = 5157.5625
$ws.Range("N17").Value = -5493.5625
$ws.Range("H43").Value = 2999.5454
$ws.Range("I43").Value = 3998.6667
$ws.Range("K43").Value = 3998.6667
$ws.Range("M43").Value = -3929.6667
$ws.Range("H70").Value = 100003
$ws.Range("J70").Value = 100003
$ws.Range("L70").Value = 300009
$ws.Range("N70").Value = -300549
$ws.Range("H73").Value = 100003
$ws.Range("J73").Value = 100003
$ws.Range("L73").Value = 300009
$ws.Range("N73").Value = -301881
$ws.Range("H76").Value = 3586.9
$ws.Range("I76").Value = 3609.375
$ws.Range("K76").Value = 3609.375
$ws.Range("M76").Value = -3294.375
$ws.Range("H79").Value = 3586.9
$ws.Range("I79").Value = 3609.375
$ws.Range("K79").Value = 3609.375
$ws.Range("M79").Value = -2517.375
$ws.Range("H86").Value = 3581
$ws.Range("I86").Value = 2940.4546
$ws.Range("K86").Value = 2940.4546
$ws.Range("M86").Value = -1817.4546
$ws.Range("H89").Value = 3581
$ws.Range("I89").Value = 2940.4546
$ws.Range("K89").Value = 14702.273
$ws.Range("M89").Value = -9086.273000000001
$ws.Range("H106").Value = 1452.7693
$ws.Range("I106").Value = 1171.7273
$ws.Range("K106").Value = 1171.7273
$ws.Range("M106").Value = -540.7273
$ws.Range("H112").Value = 1296.7142
$ws.Range("I112").Value = 1800
$ws.Range("J112").Value = 1258
$ws.Range("K112").Value = 5400
$ws.Range("L112").Value = 3774
$ws.Range("M112").Value = -4292
$ws.Range("N112").Value = -5990
$ws.Range("H137").Value = 2865.4062
$ws.Range("I137").Value = 2562.0527
$ws.Range("J137").Value = 3308.7693
$ws.Range("K137").Value = 7686.158100000001
$ws.Range("L137").Value = 9926.3079
$ws.Range("M137").Value = -5136.158100000001
$ws.Range("N137").Value = -15026.3079
$ws.Range("H141").Value = 1199.5
$ws.Range("I141").Value = 1231
$ws.Range("K141").Value = 3693
$ws.Range("M141").Value = 1487
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1505.3297
$ws.Range("I32").Value = 847.9114
$ws.Range("K32").Value = 847.9114
$ws.Range("M32").Value = -560.9114
$ws.Range("H63").Value = 7999.625
$ws.Range("I63").Value = 7599.4
$ws.Range("K63").Value = 7599.4
$ws.Range("M63").Value = -6913.4
$ws.Range("H66").Value = 7999.625
$ws.Range("I66").Value = 7599.4
$ws.Range("K66").Value = 37997
$ws.Range("M66").Value = -34565
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H88").Value = 2137
$ws.Range("I88").Value = 1799.5
$ws.Range("J88").Value = 2474.5
$ws.Range("K88").Value = 1799.5
$ws.Range("L88").Value = 2474.5
$ws.Range("M88").Value = -1393.5
$ws.Range("N88").Value = -3286.5
$ws.Range("H91").Value = 2137
$ws.Range("I91").Value = 1799.5
$ws.Range("J91").Value = 2474.5
$ws.Range("K91").Value = 1799.5
$ws.Range("L91").Value = 2474.5
$ws.Range("M91").Value = -395.5
$ws.Range("N91").Value = -5282.5
$ws.Range("H122").Value = 1353.2222
$ws.Range("I122").Value = 1272.375
$ws.Range("K122").Value = 3817.125
$ws.Range("M122").Value = -1367.125
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2000
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 2000
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H105").Value = 1386.4615
$ws.Range("I105").Value = 1335.3334
$ws.Range("K105").Value = 1335.3334
$ws.Range("M105").Value = 411.6666
$ws.Range("H107").Value = 15722.25
$ws.Range("I107").Value = 20133
$ws.Range("K107").Value = 20133
$ws.Range("M107").Value = -18213
$ws.Range("H134").Value = 1720.54
$ws.Range("I134").Value = 1084.1842
$ws.Range("J134").Value = 3735.6667
$ws.Range("K134").Value = 3252.5526
$ws.Range("L134").Value = 11207.0001
$ws.Range("M134").Value = -717.5526
$ws.Range("N134").Value = -16277.0001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 13259.25
$ws.Range("J63").Value = 13259.25
$ws.Range("L63").Value = 39777.75
$ws.Range("N63").Value = -41275.75
$ws.Range("H66").Value = 13259.25
$ws.Range("J66").Value = 13259.25
$ws.Range("L66").Value = 119333.25
$ws.Range("N66").Value = -126821.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 229.69698
$ws.Range("I2").Value = 176.85715
$ws.Range("J2").Value = 322.16666
$ws.Range("K2").Value = 176.85715
$ws.Range("L2").Value = 322.16666
$ws.Range("M2").Value = -63.85714999999999
$ws.Range("N2").Value = -548.16666
$ws.Range("H43").Value = 1753.9166
$ws.Range("I43").Value = 1753.9166
$ws.Range("K43").Value = 1753.9166
$ws.Range("M43").Value = -1602.9166
$ws.Range("H46").Value = 11393.889
$ws.Range("I46").Value = 7812.375
$ws.Range("K46").Value = 7812.375
$ws.Range("M46").Value = -7656.375
$ws.Range("H57").Value = 21833.166
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 21833.166
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 21833.166
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -23473.166
$ws.Range("H102").Value = 2618.1482
$ws.Range("I102").Value = 2603.4614
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 2603.4614
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -981.4614000000001
$ws.Range("N102").Value = -6244
$ws.Range("H123").Value = 59400
$ws.Range("I123").Value = 59000
$ws.Range("J123").Value = 59800
$ws.Range("K123").Value = 59000
$ws.Range("L123").Value = 59800
$ws.Range("M123").Value = -56550
$ws.Range("N123").Value = -64700
$ws.Range("H126").Value = 2494.8096
$ws.Range("I126").Value = 2452.6843
$ws.Range("J126").Value = 2895
$ws.Range("K126").Value = 7358.0529
$ws.Range("L126").Value = 8685
$ws.Range("M126").Value = -4888.0529
$ws.Range("N126").Value = -13625
$ws.Range("H132").Value = 30314474
$ws.Range("I132").Value = 35722000
$ws.Range("K132").Value = 107166000
$ws.Range("M132").Value = -107163470
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 445.77777
$ws.Range("I9").Value = 200.8
$ws.Range("J9").Value = 752
$ws.Range("K9").Value = 200.8
$ws.Range("L9").Value = 752
$ws.Range("M9").Value = 23.19999999999999
$ws.Range("N9").Value = -1200
$ws.Range("H93").Value = 2507.75
$ws.Range("I93").Value = 2459.3
$ws.Range("K93").Value = 2459.3
$ws.Range("M93").Value = -1211.3
$ws.Range("H98").Value = 80935.375
$ws.Range("J98").Value = 80935.375
$ws.Range("L98").Value = 80935.375
$ws.Range("N98").Value = -86925.375
$ws.Range("H101").Value = 62883.8
$ws.Range("J101").Value = 62883.8
$ws.Range("L101").Value = 62883.8
$ws.Range("N101").Value = -69373.8
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2310.0334
$ws.Range("I132").Value = 1624.0476
$ws.Range("K132").Value = 4872.142800000001
$ws.Range("M132").Value = -2342.142800000001
$ws.Range("H136").Value = 3053.3076
$ws.Range("I136").Value = 2464.2766
$ws.Range("J136").Value = 8590.200000000001
$ws.Range("K136").Value = 7392.8298
$ws.Range("L136").Value = 25770.6
$ws.Range("M136").Value = -4842.8298
$ws.Range("N136").Value = -30870.6
